# odporyspinan.xlsx - "opravena chyba pozice kurzoru v menu 1, zmenen clock na externi,
# include knihoven eaglu" - update source power resistor value and add a
# measured-vs-calculated verification table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- core change: power value used by the sheet (G2) 15000 -> 20000 ---
$ws.Range("G2").Value = 20000

# Recalculate so every dependent formula (B/H/C columns) reflects the new G2
# before we snapshot a few of the resulting values into the check table below.
$wb.Application.Calculate()

# --- fix cursor/selection position (mentioned in the commit message) ---
$ws.Application.Goto($ws.Range("V9"))

# --- add measured (Q) vs. calculated (R) verification/check table ---
$ws.Range("Q10").Value = 4.88
$ws.Range("Q11").Value = 6.93
$ws.Range("Q12").Value = 6.1
$ws.Range("Q13").Value = 9.44
$ws.Range("Q14").Value = 10.17

$ws.Range("R10").Value = $ws.Range("B2").Value2
$ws.Range("R11").Value = $ws.Range("B4").Value2
$ws.Range("R12").Value = $ws.Range("B5").Value2
$ws.Range("R13").Value = $ws.Range("B16").Value2
$ws.Range("R14").Value = $ws.Range("B14").Value2

$ws.Range("S10").Formula = "=Q10-R10"
$ws.Range("S11:S14").Formula = "=Q11-R11"

# --- conditional formatting (color scale) for the new check column,
#     same red/white/green look as the existing L4:L17 scale ---
$cf = $ws.Range("S10:S14").FormatConditions.AddColorScale(3)
$cf.ColorScaleCriteria(1).Type = 1
$cf.ColorScaleCriteria(1).FormatColor.Color = 7039480
$cf.ColorScaleCriteria(2).Type = 5
$cf.ColorScaleCriteria(2).Value = 50
$cf.ColorScaleCriteria(2).FormatColor.Color = 16776444
$cf.ColorScaleCriteria(3).Type = 2
$cf.ColorScaleCriteria(3).FormatColor.Color = 8109667
$cf.SetFirstPriority()

$wb.Save()
